$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Row 7 ("Release 1" items): drop the 4 items that move to Release 2
# (Make 911 simulation app downloadable / Create accessible link for
# parents / Create list interface... / Create readable structure...)
# Keep D7, F7, H7, T7 untouched.
# ------------------------------------------------------------------
$ws.Range("J7").Clear()
$ws.Range("N7").Clear()
$ws.Range("P7").Clear()
$ws.Range("R7").Clear()

# ------------------------------------------------------------------
# Row 9 used to hold 3 more Release-1-band items (D9, P9, R9) with the
# "item" style/height. It becomes the "Release 2" banner row instead
# (same look as the Release 1 banner row, 5).
# ------------------------------------------------------------------
$ws.Range("D9").Clear()
$ws.Range("P9").Clear()
$ws.Range("R9").Clear()

$ws.Range("B5:T5").Copy()
$ws.Range("B9:T9").PasteSpecial(-4122)
$ws.Range("B9").Value = "Release 2"
$ws.Rows.Item(9).AutoFit()

# ------------------------------------------------------------------
# Row 11 used to be the "Release 2" banner row; it becomes an item row
# holding the items that now belong to Release 2.
# ------------------------------------------------------------------
$ws.Range("B11:T11").Clear()

$ws.Range("D7").Copy()
$ws.Range("D11").PasteSpecial(-4122)
$ws.Range("D7").Copy()
$ws.Range("J11").PasteSpecial(-4122)
$ws.Range("D7").Copy()
$ws.Range("N11").PasteSpecial(-4122)
$ws.Range("D7").Copy()
$ws.Range("P11").PasteSpecial(-4122)
$ws.Range("D7").Copy()
$ws.Range("R11").PasteSpecial(-4122)
$ws.Rows.Item(11).RowHeight = 72

$ws.Range("D11").Value = "Point system is added to app to gamify it."
$ws.Range("J11").Value = "Make 911 simulation app downloadable"
$ws.Range("N11").Value = "Create accessible link for parents"
$ws.Range("P11").Value = "Create list interface to show questions and appropriate responses."
$ws.Range("R11").Value = "Display additional fire safety info in app"

# ------------------------------------------------------------------
# Row 13 already had "Link to app is accessible..." in N13; it gains
# two more Release-2 items.
# ------------------------------------------------------------------
$ws.Range("N13").Copy()
$ws.Range("P13").PasteSpecial(-4122)
$ws.Range("N13").Copy()
$ws.Range("R13").PasteSpecial(-4122)

$ws.Range("P13").Value = "Display the list of questions that 911 operators will ask"
$ws.Range("R13").Value = "Create readable structure for additional information"
